# Update the final "…" placeholder paragraph with the real changelog text,
# then add the remaining changelog paragraphs that follow it (replacing the
# trailing empty paragraph with the first new paragraph, then appending the
# rest).

$d = $word.ActiveDocument

# 1. "…" -> "Updating the document, now version 1.1.0!"
$d.Content.Find.Execute("…", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Updating the document, now version 1.1.0!", 2)

# 2. The trailing empty paragraph becomes the "Github is a version control..." paragraph.
$trailing = $d.Paragraphs.Last
$trailing.Range.InsertAfter("Github is a version control platform, used by a large number of people.")

# 3. New paragraph: "Hopefully these are enough changes, ..."
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertAfter("Hopefully these are enough changes, doesn" + [char]0x2019 + "t say anything in the rubric about marks for this.")

# 4. New paragraph with three separate runs: "Changes made by " / " I-Pudding-I" / " on 09/09/2022."
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = "<w:p $w><w:r><w:t xml:space=`"preserve`">Changes made by </w:t></w:r>" + `
       "<w:r><w:t xml:space=`"preserve`"> I-Pudding-I</w:t></w:r>" + `
       "<w:r><w:t xml:space=`"preserve`"> on 09/09/2022.</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)
